$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells F1:H1
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the header formatting (bold font, border, centered/top alignment)
# from an existing header cell so the new headers share the same style.
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# New boolean "Outliers_MAD" columns for each data row
$ws.Range("F2").Value = $false
$ws.Range("G2").Value = $false
$ws.Range("H2").Value = $false

$ws.Range("F3").Value = $true
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = $false

$ws.Range("F4").Value = $false
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = $false

$ws.Range("F5").Value = $false
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = $false

$ws.Range("F6").Value = $false
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = $false
